# Template-Penapisan.docx — "[ADD] add kbli & sector"
#
# The only content-significant change in this revision is a column-width
# re-balance on the two small tables right after the "Daftar Kegiatan" /
# "Daftar Lokasi" headings, widening room for the (new) KBLI / sector
# columns:
#
#   Table "Daftar Kegiatan" (5 columns): col3 4509dxa->3091dxa,
#                                         col4 1417dxa->2835dxa
#   Table "Daftar Lokasi"   (4 columns): col3 2545dxa->2457dxa,
#                                         col4 3881dxa->3969dxa
#
# (dxa are twentieths of a point, so Column.Width — which COM reports in
# points — is the dxa value divided by 20.)

$d = $word.ActiveDocument

# Table 2 = "Daftar Kegiatan" (Kegiatan/Nama Kegiatan/Jenis Kegiatan/...)
$tblKegiatan = $d.Tables.Item(2)
$tblKegiatan.Columns.Item(3).Width = 154.55   # 3091 dxa
$tblKegiatan.Columns.Item(4).Width = 141.75   # 2835 dxa

# Table 3 = "Daftar Lokasi" (Provinsi/Kabupaten-Kota/Alamat)
$tblLokasi = $d.Tables.Item(3)
$tblLokasi.Columns.Item(3).Width = 122.85     # 2457 dxa
$tblLokasi.Columns.Item(4).Width = 198.45     # 3969 dxa

Write-Output "Resized Daftar Kegiatan + Daftar Lokasi columns"
